$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.527.92'
$ws.Range("E2").Value = '  -0.52%  '
$ws.Range("D3").Value = '3.256.52'
$ws.Range("E3").Value = '  +3.88%  '
$ws.Range("E4").Value = '  +0.03%  '
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.65'
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = '  -0.36%  '
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.64'
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = '  +1.42%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '3.250.18'
$ws.Range("E9").Value = '  -1.26%  '
$ws.Range("E10").Value = '  -0.01%  '
$ws.Range("E11").Value = '  +1.33%  '
$ws.Range("E12").Value = '  +0.82%  '
$ws.Range("E13").Value = '  -1.75%  '
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.45'
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = '  -0.08%  '
$ws.Range("D15").Value = '3.791.30'
$ws.Range("E15").Value = '  +3.88%  '
$ws.Range("E16").Value = '  -0.20%  '
$ws.Range("D17").Value = '3.250.94'
$ws.Range("E17").Value = '  +3.99%  '
$ws.Range("D18").Value = '63.547.76'
$ws.Range("E18").Value = '  -0.43%  '
$ws.Range("E19").Value = '  +0.42%  '
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '478.47'
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = '  -0.60%  '
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.25'
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = '  -1.77%  '
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.731'
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = '  +3.84%  '
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.99'
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = '  +4.17%  '
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.85'
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = '  -4.36%  '
$ws.Range("E25").Value = '  +1.75%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("E27").Value = '  +0.26%  '
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.21'
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = '  +4.22%  '
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.10'
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = '  -0.29%  '
$ws.Range("E30").Value = '  +4.39%  '
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '27.70'
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = '  +1.99%  '
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("E33").Value = '  -3.53%  '
$ws.Range("E34").Value = '  -1.22%  '
$ws.Range("E35").Value = '  -1.01%  '
$ws.Range("E36").Value = '  -0.92%  '
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.76'
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = '  +0.39%  '
$ws.Range("D38").Value = '0.0₃0719'
$ws.Range("E38").Value = '  -2.12%  '
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0394'
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = '  +0.14%  '
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '423.44'
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = '  -0.90%  '
$ws.Range("D41").Value = '3.001.27'
$ws.Range("E41").Value = '  +4.38%  '
$ws.Range("E42").Value = '  -2.73%  '
$ws.Range("E43").Value = '  +1.17%  '
$ws.Range("E44").Value = '  -7.48%  '
$ws.Range("E45").Value = '  +2.63%  '
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.19'
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = '  +1.61%  '
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("E48").Value = '  -1.37%  '
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.98'
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = '  +1.75%  '
$ws.Range("E50").Value = '  +0.33%  '
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '122.86'
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = '  +2.15%  '
